# Update cryptos list on Sat Apr  8 16:08:43 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D include numeric-looking strings (e.g. "1.0000",
# "20.70") where trailing zeros / formatting must be preserved as text,
# so force the column to Text before writing the refreshed values.
$ws.Range("D2:D51").NumberFormat = "@"

$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '28.110.56', '  +0.26%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.873.12', '  +0.46%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.002', '  -0.22%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '313.68', '  +0.43%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.003', '  -0.01%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5107', '  +0.08%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3896', '  +0.92%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08345', '  +1.12%  '),
    @(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.118', '  +0.26%  '),
    @(11, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.229', '  -0.07%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.872.77', '  +0.79%  '),
    @(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.52', '  -0.37%  '),
    @(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.237', '  -0.14%  '),
    @(15, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.0000', '  -0.31%  '),
    @(16, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001099', '  +0.18%  '),
    @(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '90.83', '  -0.12%  '),
    @(18, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06681', '  +0.69%  '),
    @(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.73', '  -0.07%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.003', '  +0.01%  '),
    @(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.009', '  -0.44%  '),
    @(22, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '28.166.70', '  +0.38%  '),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.10', '  -0.08%  '),
    @(24, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.263', '  +1.20%  '),
    @(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '158.81', '  +1.28%  '),
    @(26, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.447', '  -3.25%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '20.70', '  +0.68%  '),
    @(28, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '126.12', '  +0.86%  '),
    @(29, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1053', '  -0.69%  '),
    @(30, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.033', '  -0.31%  '),
    @(31, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.848', '  -1.43%  '),
    @(32, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.605', '  +0.39%  '),
    @(33, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.513', '  +0.99%  '),
    @(34, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02439', '  +0.75%  '),
    @(35, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06544', '  -0.10%  '),
    @(36, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2208', '  +1.35%  '),
    @(37, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.192', '  -1.00%  '),
    @(38, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6454', '  -1.70%  '),
    @(39, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.242', '  +1.68%  '),
    @(40, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.979', '  -0.90%  '),
    @(41, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '11.28', '  +0.70%  '),
    @(42, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6072', '  -1.70%  '),
    @(43, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '13.03', '  -0.67%  '),
    @(44, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.682', '  +0.84%  '),
    @(45, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.274', '  -1.12%  '),
    @(46, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.006', '  -0.76%  '),
    @(47, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.230', '  +1.69%  '),
    @(48, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '120.55', '  +0.19%  '),
    @(49, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06893', '  +0.82%  '),
    @(50, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '77.71', '  -1.06%  '),
    @(51, 'Chiliz', 'https://coinranking.com/coin/GSCt2y6YSgO26+chiliz-chz', '0.1404', '  +8.98%  ')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
